$d = $word.ActiveDocument

# Target: the empty paragraph right after "Trần Hưng Thuận" (under "Tác giả trả lời")
# and right before "Võ Xuân Tiến". It is paragraph index 21.
$target = $d.Paragraphs(21)

# 1) Fill in the date/time line "(15h, 19/3)" on the currently-empty paragraph.
$target.Range.Text = "(15h, 19/3)"

# 2) Insert a new paragraph after it for "Đã sửa", with ListParagraph style and
#    bullet numbering (ilvl 0, numId 8) matching the rest of this person's list.
$target.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(22)
$p2.Range.Text = "Đã sửa"
$p2.Range.ListFormat.ApplyListTemplateWithLevel($d.ListTemplates.Item(8), $false, 1, $false, 0)

# 3) Insert a third paragraph after that with the longer explanation text, same numbering.
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(23)
$p3.Range.Text = "Theo ý kiến bàn bạc thì có ý kiến thêm nhân viên, nhưng lúc thiết kế testcase đã bỏ đi actor nhân viên, nên chỉ còn 2 actor thôi. Đã chỉnh sửa “Khách hàng” thành “User” để làm rõ thêm"
$p3.Range.ListFormat.ApplyListTemplateWithLevel($d.ListTemplates.Item(8), $false, 1, $false, 0)
